$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 changes
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 2.25
$ws.Range("J3").Value = 3.6
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 3
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.85
$ws.Range("W3").Value = 10
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 34
$ws.Range("AA3").Value = 23
$ws.Range("AC3").Value = 10
$ws.Range("AF3").Value = 41
$ws.Range("AI3").Value = 11
$ws.Range("AJ3").Value = 9.5
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 19
$ws.Range("AM3").Value = 29
$ws.Range("AN3").Value = 5
$ws.Range("AO3").Value = 17
$ws.Range("AQ3").Value = 51
$ws.Range("AS3").Value = 151
$ws.Range("AW3").Value = 4.33
$ws.Range("AX3").Value = 13
$ws.Range("AY3").Value = 23

# Row 4 changes
$ws.Range("G4").Value = 1.67
$ws.Range("I4").Value = 5.25
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("Z4").Value = 12
$ws.Range("AI4").Value = 26
$ws.Range("AQ4").Value = 29
$ws.Range("AT4").Value = 2.63
$ws.Range("AZ4").Value = 126
